$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 269.5
$ws.Range("I5").Value = 229
$ws.Range("K5").Value = 229
$ws.Range("M5").Value = -114
$ws.Range("H40").Value = 33433.145
$ws.Range("I40").Value = 22806.4
$ws.Range("K40").Value = 22806.4
$ws.Range("M40").Value = -22631.4
$ws.Range("H43").Value = 3856904
$ws.Range("I43").Value = 3856904
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3856904
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3856835
$ws.Range("N43").Value = ""
$ws.Range("H62").Value = 4168
$ws.Range("I62").Value = 3639.1765
$ws.Range("K62").Value = 3639.1765
$ws.Range("M62").Value = -3015.1765
$ws.Range("H64").Value = 2178486.8
$ws.Range("I64").Value = 6215114.5
$ws.Range("K64").Value = 6215114.5
$ws.Range("M64").Value = -6214866.5
$ws.Range("H65").Value = 4168
$ws.Range("I65").Value = 3639.1765
$ws.Range("K65").Value = 18195.8825
$ws.Range("M65").Value = -15075.8825
$ws.Range("H67").Value = 2178486.8
$ws.Range("I67").Value = 6215114.5
$ws.Range("K67").Value = 6215114.5
$ws.Range("M67").Value = -6214256.5
$ws.Range("H98").Value = 1815.2858
$ws.Range("I98").Value = 1246.2106
$ws.Range("K98").Value = 1246.2106
$ws.Range("M98").Value = 251.7893999999999
$ws.Range("H111").Value = 7370.75
$ws.Range("I111").Value = 7709.4287
$ws.Range("K111").Value = 23128.2861
$ws.Range("M111").Value = -20061.2861
$ws.Range("H112").Value = 3096.6086
$ws.Range("J112").Value = 3115.4443
$ws.Range("L112").Value = 9346.332900000001
$ws.Range("N112").Value = -11562.3329
$ws.Range("H116").Value = 21218160
$ws.Range("I116").Value = 27090454
$ws.Range("J116").Value = 17862564
$ws.Range("K116").Value = 27090454
$ws.Range("L116").Value = 17862564
$ws.Range("M116").Value = -27087012
$ws.Range("N116").Value = -17869448
$ws.Range("H122").Value = 1815.2858
$ws.Range("I122").Value = 1246.2106
$ws.Range("K122").Value = 3738.6318
$ws.Range("M122").Value = -1288.6318
$ws.Range("H132").Value = 135425.64
$ws.Range("I132").Value = 533812.5
$ws.Range("K132").Value = 1601437.5
$ws.Range("M132").Value = -1598907.5
$ws.Range("H135").Value = 5984.5
$ws.Range("I135").Value = 907.2308
$ws.Range("K135").Value = 8165.077200000001
$ws.Range("M135").Value = -5630.077200000001
$ws.Range("H137").Value = 19610994
$ws.Range("I137").Value = 1586.7273
$ws.Range("J137").Value = 55561576
$ws.Range("K137").Value = 4760.1819
$ws.Range("L137").Value = 166684728
$ws.Range("M137").Value = -2210.1819
$ws.Range("N137").Value = -166689828
$ws.Range("H138").Value = 5674.3784
$ws.Range("I138").Value = 2319.389
$ws.Range("J138").Value = 6752.768
$ws.Range("K138").Value = 6958.167
$ws.Range("L138").Value = 20258.304
$ws.Range("M138").Value = -1818.167
$ws.Range("N138").Value = -30538.304
$ws.Range("H141").Value = 4358.0625
$ws.Range("I141").Value = 3909.2144
$ws.Range("K141").Value = 11727.6432
$ws.Range("M141").Value = -6547.643199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7010.7427
$ws.Range("I32").Value = 4553.1816
$ws.Range("K32").Value = 4553.1816
$ws.Range("M32").Value = -4266.1816
$ws.Range("H102").Value = 806711.5
$ws.Range("I102").Value = 979400
$ws.Range("J102").Value = 831.6667
$ws.Range("K102").Value = 979400
$ws.Range("L102").Value = 831.6667
$ws.Range("M102").Value = -977778
$ws.Range("N102").Value = -4075.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 244.5
$ws.Range("I22").Value = 236.57143
$ws.Range("K22").Value = 236.57143
$ws.Range("M22").Value = 113.42857
$ws.Range("H31").Value = 4263.1406
$ws.Range("I31").Value = 2910.2083
$ws.Range("J31").Value = 7086.6523
$ws.Range("K31").Value = 2910.2083
$ws.Range("L31").Value = 7086.6523
$ws.Range("M31").Value = -2615.2083
$ws.Range("N31").Value = -7676.6523
$ws.Range("H34").Value = 4263.1406
$ws.Range("I34").Value = 2910.2083
$ws.Range("J34").Value = 7086.6523
$ws.Range("K34").Value = 2910.2083
$ws.Range("L34").Value = 7086.6523
$ws.Range("M34").Value = -2708.2083
$ws.Range("N34").Value = -7490.6523
$ws.Range("H58").Value = 627753.25
$ws.Range("I58").Value = 1002517
$ws.Range("K58").Value = 1002517
$ws.Range("M58").Value = -1002314
$ws.Range("H122").Value = 4812033
$ws.Range("I122").Value = 10994898
$ws.Range("K122").Value = 32984694
$ws.Range("M122").Value = -32982244
$ws.Range("H132").Value = 2212.3
$ws.Range("I132").Value = 1020.8333
$ws.Range("K132").Value = 3062.4999
$ws.Range("M132").Value = -532.4998999999998
$ws.Range("H136").Value = 627753.25
$ws.Range("I136").Value = 1002517
$ws.Range("K136").Value = 3007551
$ws.Range("M136").Value = -3005001
$ws.Range("H139").Value = 77500
$ws.Range("J139").Value = 105000
$ws.Range("L139").Value = 105000
$ws.Range("N139").Value = -115280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47397748
$ws.Range("I4").Value = 1235543.4
$ws.Range("K4").Value = 3706630.2
$ws.Range("M4").Value = -3706518.2
$ws.Range("H50").Value = 1225.5
$ws.Range("I50").Value = 445
$ws.Range("K50").Value = 1335
$ws.Range("M50").Value = -854
$ws.Range("H53").Value = 1225.5
$ws.Range("I53").Value = 445
$ws.Range("K53").Value = 1335
$ws.Range("M53").Value = -854

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1650130.2
$ws.Range("I70").Value = 2652976
$ws.Range("J70").Value = 9109.727999999999
$ws.Range("K70").Value = 2652976
$ws.Range("L70").Value = 9109.727999999999
$ws.Range("M70").Value = -2652706
$ws.Range("N70").Value = -9649.727999999999
$ws.Range("H73").Value = 1650130.2
$ws.Range("I73").Value = 2652976
$ws.Range("J73").Value = 9109.727999999999
$ws.Range("K73").Value = 2652976
$ws.Range("L73").Value = 9109.727999999999
$ws.Range("M73").Value = -2652040
$ws.Range("N73").Value = -10981.728
$ws.Range("H80").Value = 42134.58
$ws.Range("I80").Value = 79002.766
$ws.Range("J80").Value = 5266.385
$ws.Range("K80").Value = 79002.766
$ws.Range("L80").Value = 5266.385
$ws.Range("M80").Value = -78004.766
$ws.Range("N80").Value = -7262.385
$ws.Range("H83").Value = 42134.58
$ws.Range("I83").Value = 79002.766
$ws.Range("J83").Value = 5266.385
$ws.Range("K83").Value = 395013.83
$ws.Range("L83").Value = 26331.925
$ws.Range("M83").Value = -390021.83
$ws.Range("N83").Value = -36315.925
$ws.Range("H132").Value = 6733.64
$ws.Range("I132").Value = 6058
$ws.Range("J132").Value = 8169.375
$ws.Range("K132").Value = 18174
$ws.Range("L132").Value = 24508.125
$ws.Range("M132").Value = -15644
$ws.Range("N132").Value = -29568.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4882.268
$ws.Range("I7").Value = 3917.742
$ws.Range("K7").Value = 3917.742
$ws.Range("M7").Value = -3805.742
$ws.Range("H16").Value = 2243.353
$ws.Range("I16").Value = 1509.6428
$ws.Range("J16").Value = 5667.3335
$ws.Range("K16").Value = 1509.6428
$ws.Range("L16").Value = 5667.3335
$ws.Range("M16").Value = -1339.6428
$ws.Range("N16").Value = -6007.3335
$ws.Range("H22").Value = 1077.4615
$ws.Range("I22").Value = 850.36365
$ws.Range("J22").Value = 1244
$ws.Range("K22").Value = 850.36365
$ws.Range("L22").Value = 1244
$ws.Range("M22").Value = -555.36365
$ws.Range("N22").Value = -1834
$ws.Range("H27").Value = 1077.4615
$ws.Range("I27").Value = 850.36365
$ws.Range("J27").Value = 1244
$ws.Range("K27").Value = 850.36365
$ws.Range("L27").Value = 1244
$ws.Range("M27").Value = -743.36365
$ws.Range("N27").Value = -1458
$ws.Range("H46").Value = 5390.4614
$ws.Range("I46").Value = 2575
$ws.Range("J46").Value = 6235.1
$ws.Range("K46").Value = 2575
$ws.Range("L46").Value = 6235.1
$ws.Range("M46").Value = -2387
$ws.Range("N46").Value = -6611.1
$ws.Range("H55").Value = 466.4
$ws.Range("J55").Value = 410
$ws.Range("L55").Value = 410
$ws.Range("N55").Value = -756
$ws.Range("H126").Value = 4882.268
$ws.Range("I126").Value = 3917.742
$ws.Range("K126").Value = 11753.226
$ws.Range("M126").Value = -9283.226000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 37061308
$ws.Range("I132").Value = 4301122
$ws.Range("J132").Value = 250002500
$ws.Range("K132").Value = 12903366
$ws.Range("L132").Value = 750007500
$ws.Range("M132").Value = -12900836
$ws.Range("N132").Value = -750012560
